$d = $word.ActiveDocument

# The last paragraph in the document is the "Silly Love Songs" list item,
# which also carries the trailing _GoBack bookmark. We want three new list
# entries to end up AFTER it, with the bookmark ending up on the very last
# paragraph ("Ed Sheeran Album").
#
# Rather than fighting bookmark relocation, we keep this paragraph node
# (and its bookmark) in place structurally, insert three fresh paragraphs
# BEFORE it (inheriting the same list formatting), and then shuffle the
# text: the three new paragraphs get "Silly Love Songs...", "Love
# Yourself...", "Moon River...", while the original (bookmarked) paragraph
# becomes "Ed Sheeran Album" - ending up last, bookmark and all.

$anchorIndex = $d.Paragraphs.Count
$anchor = $d.Paragraphs($anchorIndex)

# Insert three blank paragraphs immediately before the anchor paragraph;
# they inherit the ListParagraph style + numbering from the anchor.
$anchor.Range.InsertParagraphBefore()
$anchor.Range.InsertParagraphBefore()
$anchor.Range.InsertParagraphBefore()

# After the inserts, paragraph indices are:
#   anchorIndex    -> new blank paragraph #1
#   anchorIndex+1  -> new blank paragraph #2
#   anchorIndex+2  -> new blank paragraph #3
#   anchorIndex+3  -> original anchor paragraph (still holds the bookmark)
$d.Paragraphs($anchorIndex).Range.Text     = "Silly Love Songs – Paul McCartney"
$d.Paragraphs($anchorIndex + 1).Range.Text = "Love Yourself – Justin Bieber"
$d.Paragraphs($anchorIndex + 2).Range.Text = "Moon River – Audrey Hepburn"
$d.Paragraphs($anchorIndex + 3).Range.Text = "Ed Sheeran Album"

Write-Host "Inserted 3 new entries"
